$wb = $excel.ActiveWorkbook

# ---- Sheet: BSM ----
$wsBSM = $wb.Worksheets.Item("BSM")
$wsBSM.Range("H117").Value = 40742
$wsBSM.Range("I117").Value = 0
$wsBSM.Range("J117").Value = 40742
$wsBSM.Range("K117").Value = 0
$wsBSM.Range("L117").Value = 40742
$wsBSM.Range("N117").Value = -49920
$wsBSM.Range("H118").Value = 78333.336
$wsBSM.Range("I118").Value = 0
$wsBSM.Range("J118").Value = 78333.336
$wsBSM.Range("K118").Value = 0
$wsBSM.Range("L118").Value = 78333.336
$wsBSM.Range("N118").Value = -81647.336
$wsBSM.Range("H119").Value = 39508.8
$wsBSM.Range("I119").Value = 0
$wsBSM.Range("J119").Value = 39508.8
$wsBSM.Range("K119").Value = 0
$wsBSM.Range("L119").Value = 39508.8
$wsBSM.Range("N119").Value = -49184.8
$wsBSM.Range("H120").Value = 40761
$wsBSM.Range("I120").Value = 0
$wsBSM.Range("J120").Value = 40761
$wsBSM.Range("K120").Value = 0
$wsBSM.Range("L120").Value = 40761
$wsBSM.Range("N120").Value = -50437
$wsBSM.Range("H122").Value = 50780
$wsBSM.Range("I122").Value = 0
$wsBSM.Range("J122").Value = 50780
$wsBSM.Range("K122").Value = 0
$wsBSM.Range("L122").Value = 50780
$wsBSM.Range("N122").Value = -60580
$wsBSM.Range("H123").Value = 41996.668
$wsBSM.Range("I123").Value = 0
$wsBSM.Range("J123").Value = 41996.668
$wsBSM.Range("K123").Value = 0
$wsBSM.Range("L123").Value = 41996.668
$wsBSM.Range("N123").Value = -51796.668
$wsBSM.Range("H124").Value = 37593.332
$wsBSM.Range("I124").Value = 0
$wsBSM.Range("J124").Value = 37593.332
$wsBSM.Range("K124").Value = 0
$wsBSM.Range("L124").Value = 37593.332
$wsBSM.Range("N124").Value = -47413.332
$wsBSM.Range("H125").Value = 40780
$wsBSM.Range("I125").Value = 0
$wsBSM.Range("J125").Value = 40780
$wsBSM.Range("K125").Value = 0
$wsBSM.Range("L125").Value = 40780
$wsBSM.Range("N125").Value = -50620
$wsBSM.Range("H126").Value = 0
$wsBSM.Range("I126").Value = 0
$wsBSM.Range("J126").Value = 0
$wsBSM.Range("K126").Value = 0
$wsBSM.Range("L126").Value = 0
$wsBSM.Range("H127").Value = 50000
$wsBSM.Range("I127").Value = 0
$wsBSM.Range("J127").Value = 50000
$wsBSM.Range("K127").Value = 0
$wsBSM.Range("L127").Value = 50000
$wsBSM.Range("N127").Value = -59920
$wsBSM.Range("H128").Value = 1000
$wsBSM.Range("I128").Value = 1000
$wsBSM.Range("J128").Value = 0
$wsBSM.Range("K128").Value = 3000
$wsBSM.Range("L128").Value = 0
$wsBSM.Range("M128").Value = -510
$wsBSM.Range("H129").Value = 49199.6
$wsBSM.Range("I129").Value = 0
$wsBSM.Range("J129").Value = 49199.6
$wsBSM.Range("K129").Value = 0
$wsBSM.Range("L129").Value = 49199.6
$wsBSM.Range("N129").Value = -59199.6
$wsBSM.Range("H130").Value = 40624
$wsBSM.Range("I130").Value = 0
$wsBSM.Range("J130").Value = 40624
$wsBSM.Range("K130").Value = 0
$wsBSM.Range("L130").Value = 40624
$wsBSM.Range("N130").Value = -50664
$wsBSM.Range("H131").Value = 27819.334
$wsBSM.Range("I131").Value = 0
$wsBSM.Range("J131").Value = 27819.334
$wsBSM.Range("K131").Value = 0
$wsBSM.Range("L131").Value = 27819.334
$wsBSM.Range("N131").Value = -37899.334
$wsBSM.Range("H132").Value = 40125
$wsBSM.Range("I132").Value = 0
$wsBSM.Range("J132").Value = 40125
$wsBSM.Range("K132").Value = 0
$wsBSM.Range("L132").Value = 40125
$wsBSM.Range("N132").Value = -50245
$wsBSM.Range("H133").Value = 29800
$wsBSM.Range("I133").Value = 0
$wsBSM.Range("J133").Value = 29800
$wsBSM.Range("K133").Value = 0
$wsBSM.Range("L133").Value = 29800
$wsBSM.Range("N133").Value = -39920
$wsBSM.Range("H134").Value = 2168.5
$wsBSM.Range("I134").Value = 2224.111
$wsBSM.Range("J134").Value = 2043.375
$wsBSM.Range("K134").Value = 6672.333
$wsBSM.Range("L134").Value = 6130.125
$wsBSM.Range("M134").Value = -4137.333
$wsBSM.Range("N134").Value = -11200.125
$wsBSM.Range("H135").Value = 32000
$wsBSM.Range("I135").Value = 0
$wsBSM.Range("J135").Value = 32000
$wsBSM.Range("K135").Value = 0
$wsBSM.Range("L135").Value = 32000
$wsBSM.Range("N135").Value = -42140
$wsBSM.Range("H137").Value = 0
$wsBSM.Range("I137").Value = 0
$wsBSM.Range("J137").Value = 0
$wsBSM.Range("K137").Value = 0
$wsBSM.Range("L137").Value = 0
$wsBSM.Range("H138").Value = 30000
$wsBSM.Range("I138").Value = 0
$wsBSM.Range("J138").Value = 30000
$wsBSM.Range("K138").Value = 0
$wsBSM.Range("L138").Value = 30000
$wsBSM.Range("N138").Value = -40280
$wsBSM.Range("H139").Value = 50000
$wsBSM.Range("I139").Value = 0
$wsBSM.Range("J139").Value = 50000
$wsBSM.Range("K139").Value = 0
$wsBSM.Range("L139").Value = 50000
$wsBSM.Range("N139").Value = -60280
$wsBSM.Range("H140").Value = 280000
$wsBSM.Range("I140").Value = 0
$wsBSM.Range("J140").Value = 280000
$wsBSM.Range("K140").Value = 0
$wsBSM.Range("L140").Value = 280000
$wsBSM.Range("N140").Value = -290360
$wsBSM.Range("H141").Value = 0
$wsBSM.Range("I141").Value = 0
$wsBSM.Range("J141").Value = 0
$wsBSM.Range("K141").Value = 0
$wsBSM.Range("L141").Value = 0

# ---- Sheet: CUL ----
$wsCUL = $wb.Worksheets.Item("CUL")
$wsCUL.Range("H131").Value = 971.4375
$wsCUL.Range("I131").Value = 478.33334
$wsCUL.Range("K131").Value = 1435.00002
$wsCUL.Range("M131").Value = 3604.99998

# ---- Sheet: WVR ----
$wsWVR = $wb.Worksheets.Item("WVR")
$wsWVR.Range("H119").Value = 0
$wsWVR.Range("I119").Value = 0
$wsWVR.Range("J119").Value = 0
$wsWVR.Range("K119").Value = 0
$wsWVR.Range("L119").Value = 0
$wsWVR.Range("H120").Value = 40000
$wsWVR.Range("I120").Value = 0
$wsWVR.Range("J120").Value = 40000
$wsWVR.Range("K120").Value = 0
$wsWVR.Range("L120").Value = 40000
$wsWVR.Range("N120").Value = -49676
$wsWVR.Range("H121").Value = 50000
$wsWVR.Range("I121").Value = 0
$wsWVR.Range("J121").Value = 50000
$wsWVR.Range("K121").Value = 0
$wsWVR.Range("L121").Value = 50000
$wsWVR.Range("N121").Value = -53494
$wsWVR.Range("H122").Value = 2212.2222
$wsWVR.Range("I122").Value = 2273
$wsWVR.Range("J122").Value = 1999.5
$wsWVR.Range("K122").Value = 6819
$wsWVR.Range("L122").Value = 5998.5
$wsWVR.Range("M122").Value = -4369
$wsWVR.Range("N122").Value = -10898.5
$wsWVR.Range("H123").Value = 0
$wsWVR.Range("I123").Value = 0
$wsWVR.Range("J123").Value = 0
$wsWVR.Range("K123").Value = 0
$wsWVR.Range("L123").Value = 0
$wsWVR.Range("H124").Value = 15000
$wsWVR.Range("I124").Value = 0
$wsWVR.Range("J124").Value = 15000
$wsWVR.Range("K124").Value = 0
$wsWVR.Range("L124").Value = 15000
$wsWVR.Range("N124").Value = -24820
$wsWVR.Range("H125").Value = 50000
$wsWVR.Range("I125").Value = 0
$wsWVR.Range("J125").Value = 50000
$wsWVR.Range("K125").Value = 0
$wsWVR.Range("L125").Value = 50000
$wsWVR.Range("N125").Value = -59840
$wsWVR.Range("H126").Value = 991.62964
$wsWVR.Range("I126").Value = 950.5599999999999
$wsWVR.Range("J126").Value = 1505
$wsWVR.Range("K126").Value = 2851.68
$wsWVR.Range("L126").Value = 4515
$wsWVR.Range("M126").Value = -381.6799999999998
$wsWVR.Range("N126").Value = -9455
$wsWVR.Range("H127").Value = 50000
$wsWVR.Range("I127").Value = 0
$wsWVR.Range("J127").Value = 50000
$wsWVR.Range("K127").Value = 0
$wsWVR.Range("L127").Value = 50000
$wsWVR.Range("N127").Value = -59920
$wsWVR.Range("H128").Value = 36505
$wsWVR.Range("I128").Value = 0
$wsWVR.Range("J128").Value = 36505
$wsWVR.Range("K128").Value = 0
$wsWVR.Range("L128").Value = 36505
$wsWVR.Range("N128").Value = -46465
$wsWVR.Range("H129").Value = 49214.5
$wsWVR.Range("I129").Value = 0
$wsWVR.Range("J129").Value = 49214.5
$wsWVR.Range("K129").Value = 0
$wsWVR.Range("L129").Value = 49214.5
$wsWVR.Range("N129").Value = -59214.5
$wsWVR.Range("H130").Value = 19000
$wsWVR.Range("I130").Value = 0
$wsWVR.Range("J130").Value = 19000
$wsWVR.Range("K130").Value = 0
$wsWVR.Range("L130").Value = 19000
$wsWVR.Range("N130").Value = -29040
$wsWVR.Range("H131").Value = 50000
$wsWVR.Range("I131").Value = 0
$wsWVR.Range("J131").Value = 50000
$wsWVR.Range("K131").Value = 0
$wsWVR.Range("L131").Value = 50000
$wsWVR.Range("N131").Value = -60080
$wsWVR.Range("H132").Value = 1952.4615
$wsWVR.Range("I132").Value = 1556.6
$wsWVR.Range("J132").Value = 2199.875
$wsWVR.Range("K132").Value = 4669.799999999999
$wsWVR.Range("L132").Value = 6599.625
$wsWVR.Range("M132").Value = -2139.799999999999
$wsWVR.Range("N132").Value = -11659.625
$wsWVR.Range("H133").Value = 29200
$wsWVR.Range("I133").Value = 0
$wsWVR.Range("J133").Value = 29200
$wsWVR.Range("K133").Value = 0
$wsWVR.Range("L133").Value = 29200
$wsWVR.Range("N133").Value = -39320
$wsWVR.Range("H135").Value = 35000
$wsWVR.Range("I135").Value = 0
$wsWVR.Range("J135").Value = 35000
$wsWVR.Range("K135").Value = 0
$wsWVR.Range("L135").Value = 35000
$wsWVR.Range("N135").Value = -45140
$wsWVR.Range("H136").Value = 2820.0679
$wsWVR.Range("I136").Value = 922.2059
$wsWVR.Range("J136").Value = 5401.16
$wsWVR.Range("K136").Value = 2766.6177
$wsWVR.Range("L136").Value = 16203.48
$wsWVR.Range("M136").Value = -216.6177000000002
$wsWVR.Range("N136").Value = -21303.48
$wsWVR.Range("H137").Value = 0
$wsWVR.Range("I137").Value = 0
$wsWVR.Range("J137").Value = 0
$wsWVR.Range("K137").Value = 0
$wsWVR.Range("L137").Value = 0
$wsWVR.Range("H138").Value = 50000
$wsWVR.Range("I138").Value = 0
$wsWVR.Range("J138").Value = 50000
$wsWVR.Range("K138").Value = 0
$wsWVR.Range("L138").Value = 50000
$wsWVR.Range("N138").Value = -60280
$wsWVR.Range("H139").Value = 0
$wsWVR.Range("I139").Value = 0
$wsWVR.Range("J139").Value = 0
$wsWVR.Range("K139").Value = 0
$wsWVR.Range("L139").Value = 0
$wsWVR.Range("H140").Value = 38561.5
$wsWVR.Range("I140").Value = 0
$wsWVR.Range("J140").Value = 38561.5
$wsWVR.Range("K140").Value = 0
$wsWVR.Range("L140").Value = 38561.5
$wsWVR.Range("N140").Value = -48921.5
$wsWVR.Range("H141").Value = 45000
$wsWVR.Range("I141").Value = 0
$wsWVR.Range("J141").Value = 45000
$wsWVR.Range("K141").Value = 0
$wsWVR.Range("L141").Value = 45000
$wsWVR.Range("N141").Value = -55360
